$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data
$ws.Range("D2").Value = "24.339.18"
$ws.Range("E2").Value = "  +9.48%  "
$ws.Range("D3").Value = "1.673.57"
$ws.Range("E3").Value = "  +4.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.36"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3678"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3416"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.50"
$ws.Range("E9").Value = "  +14.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.154"
$ws.Range("E10").Value = "  +3.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07192"
$ws.Range("E11").Value = "  +3.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.114"
$ws.Range("E13").Value = "  +4.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.98"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.698"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "1.674.71"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06642"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "79.99"
$ws.Range("E20").Value = "  +4.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.39"
$ws.Range("E21").Value = "  +3.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.077"
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.16"
$ws.Range("E23").Value = "  +4.30%  "
$ws.Range("D24").Value = "24.325.93"
$ws.Range("E24").Value = "  +9.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.438"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.632"
$ws.Range("E26").Value = "  +4.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.48"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.34"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").Value = "1.861.87"
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.49"
$ws.Range("E30").Value = "  +4.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.238"
$ws.Range("E31").Value = "  +5.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.049"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9629"
$ws.Range("E33").Value = "  +4.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08441"
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.676"
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.23"
$ws.Range("E36").Value = "  +4.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06363"
$ws.Range("E37").Value = "  +5.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.265"
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02307"
$ws.Range("E39").Value = "  +5.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.631"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.236"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2075"
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6045"
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.740"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5847"
$ws.Range("E47").Value = "  +4.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.37"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.004"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07132"
$ws.Range("E50").Value = "  +5.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.49"
$ws.Range("E51").Value = "  +3.86%  "
